$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting existing rows 47-100 down to 48-101
# (dimension grows from A1:R100 to A1:R101).
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with its data. The row carries over the
# same market/category metadata as the record that used to occupy row 47
# (now at row 48), but with a new date and new volume/price figures.
$ws.Cells.Item(47, 1).Value = 11
$ws.Cells.Item(47, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(47, 3).Value = "Bíobío"
$ws.Cells.Item(47, 4).Value = 44757
$ws.Cells.Item(47, 5).Value = 8
$ws.Cells.Item(47, 6).Value = 100112001
$ws.Cells.Item(47, 7).Value = "Berenjena"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 100
$ws.Cells.Item(47, 11).Value = 12000
$ws.Cells.Item(47, 12).Value = 13000
$ws.Cells.Item(47, 13).Value = 12500
$ws.Cells.Item(47, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 208
$ws.Cells.Item(47, 17).Value = 60
$ws.Cells.Item(47, 18).Value = "Hortaliza"
